# Fungal_Models.xlsx — add the two new Aspergillus niger (iJB1325) SBML rows
# and widen the Species column, per the "adding Aspergillus niger SBMLs and
# updating Fungal_Models.xlsx" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Duplicate row 20's formatting down into two new rows (21 & 22) ---
# Row 20 only has cells in columns A, C, D, F (B/E are blank there), but the
# row itself carries a "customFormat" row style that blank B20/E20 cells
# would inherit. Populate them momentarily so the copy below also carries
# that formatting into the new rows, then restore row 20 to its original
# (cell-less) state.
$ws.Range("B20").Value = 1
$ws.Range("E20").Value = "x"

$ws.Rows("20:20").Copy()
$ws.Rows("21:21").Insert()
$ws.Rows("20:20").Copy()
$ws.Rows("22:22").Insert()

$ws.Range("B20").Clear()
$ws.Range("E20").Clear()

# B21/B22 (PubMed ID) use the hyperlink-style formatting (same as column F),
# not the plain formatting row 20 supplied - copy that format over.
$ws.Range("F5").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B22").PasteSpecial(-4122)

# --- 2. Attach the hyperlinks first ---
# (Hyperlinks.Add's TextToDisplay writes text into the cell; doing this
# before the final Value assignments below means the real numeric / string
# cell content wins, while the hyperlink's `display` attribute still ends
# up correct.)
$ws.Hyperlinks.Add($ws.Range("B21"), "https://www.ncbi.nlm.nih.gov/pubmed/30275963", [Type]::Missing, [Type]::Missing, "https://www.ncbi.nlm.nih.gov/pubmed/30275963")
$ws.Hyperlinks.Add($ws.Range("B22"), "https://www.ncbi.nlm.nih.gov/pubmed/30275963", [Type]::Missing, [Type]::Missing, "https://www.ncbi.nlm.nih.gov/pubmed/30275963")
$ws.Hyperlinks.Add($ws.Range("F21"), "https://doi.org/10.1186/s40694-018-0060-7", [Type]::Missing, [Type]::Missing, "BMC (SBML)")
$ws.Hyperlinks.Add($ws.Range("F22"), "https://doi.org/10.1186/s40694-018-0060-7", [Type]::Missing, [Type]::Missing, "BMC (SBML)")

# --- 3. Fill in the new row content ---
# Row 21: Aspergillus niger ATCC 1015
$ws.Range("A21").Value = "iJB1325"
$ws.Range("B21").Value = 30275963
$ws.Range("C21").Value = 2018
$ws.Range("D21").Value = "Aspergillus niger ATCC 1015"
$ws.Range("E21").Value = "A community-driven reconstruction of the Aspergillus niger metabolic network."
$ws.Range("F21").Value = "BMC (SBML)"

# Row 22: Aspergillus niger CBS 513.88
$ws.Range("A22").Value = "iJB1325"
$ws.Range("B22").Value = 30275963
$ws.Range("C22").Value = 2018
$ws.Range("D22").Value = "Aspergillus niger CBS 513.88"
$ws.Range("E22").Value = "A community-driven reconstruction of the Aspergillus niger metabolic network."
$ws.Range("F22").Value = "BMC (SBML)"

# --- 4. Widen the Species column (D) to fit the new, longer entries ---
$ws.Columns("D").ColumnWidth = 21.42

# --- 5. Leave the selection where the author ended up ---
$ws.Range("B22").Select()
